# The commit adds one new daily price-report record for "Zapallo italiano"
# (Vega Central Mapocho de Santiago) to the weekly consolidated sheet.
# The new record is inserted as row 277, which pushes every existing row
# from 277..365 down by one (to 278..366) and grows the used range from
# A1:R365 to A1:R366.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 277 - this shifts rows
# 277..365 down to 278..366 and extends the sheet dimension accordingly.
$ws.Rows.Item(277).Insert()

# Populate the newly inserted row 277 with the new observation.
$ws.Range("A277").Value = 9
$ws.Range("B277").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C277").Value = "Metropolitana"
$ws.Range("D277").Value = 44663
$ws.Range("E277").Value = 13
$ws.Range("F277").Value = 100112032
$ws.Range("G277").Value = "Zapallo italiano"
$ws.Range("H277").Value = "Sin especificar"
$ws.Range("I277").Value = "Primera"
$ws.Range("J277").Value = 79
$ws.Range("K277").Value = 9000
$ws.Range("L277").Value = 10000
$ws.Range("M277").Value = 9494
$ws.Range("N277").Value = "`$/caja 60 unidades"
$ws.Range("O277").Value = "Región Metropolitana"
$ws.Range("P277").Value = 158
$ws.Range("Q277").Value = 60
$ws.Range("R277").Value = "Hortaliza"
